# Apply "ajuste: corrigindo as categorias" edit:
#  - Add two new columns: S = "Idade ignorada", T = "Total"
#  - Add two new rows: 7 = "Outros", 8 = "Total" (grand total row)
#  - Existing rows 2-6 keep their text/values, gaining blank S cells and a new T (row total)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Cells.Item(1, 19).Value = "Idade ignorada"   # S1
$ws.Cells.Item(1, 20).Value = "Total"            # T1

# --- Row totals for existing rows (column T), and blank placeholders for column S ---
# Row 2: Doenças do aparelho circulatório
$ws.Range("S2").Borders.LineStyle = -4142
$ws.Cells.Item(2, 20).Value = 1997

# Row 3: Doenças do aparelho geniturinário
$ws.Range("S3").Borders.LineStyle = -4142
$ws.Cells.Item(3, 20).Value = 279

# Row 4: Doenças do aparelho respiratório
$ws.Range("S4").Borders.LineStyle = -4142
$ws.Cells.Item(4, 20).Value = 1122

# Row 5: Doenças endócrinas, nutricionais e metabólicas
$ws.Range("S5").Borders.LineStyle = -4142
$ws.Cells.Item(5, 20).Value = 282

# Row 6: Neoplasmas
$ws.Range("S6").Borders.LineStyle = -4142
$ws.Cells.Item(6, 20).Value = 1420

# --- New row 7: Outros ---
$ws.Cells.Item(7, 1).Value = "Outros"
$ws.Cells.Item(7, 2).Value = 144
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = 20
$ws.Cells.Item(7, 5).Value = 55
$ws.Cells.Item(7, 6).Value = 64
$ws.Cells.Item(7, 7).Value = 80
$ws.Cells.Item(7, 8).Value = 103
$ws.Cells.Item(7, 9).Value = 97
$ws.Cells.Item(7, 10).Value = 95
$ws.Cells.Item(7, 11).Value = 121
$ws.Cells.Item(7, 12).Value = 113
$ws.Cells.Item(7, 13).Value = 120
$ws.Cells.Item(7, 14).Value = 107
$ws.Cells.Item(7, 15).Value = 113
$ws.Cells.Item(7, 16).Value = 107
$ws.Cells.Item(7, 17).Value = 137
$ws.Cells.Item(7, 18).Value = 483
$ws.Cells.Item(7, 19).Value = 1
$ws.Cells.Item(7, 20).Value = 1963

# --- New row 8: Total (grand total) ---
$ws.Cells.Item(8, 1).Value = "Total"
$ws.Cells.Item(8, 2).Value = 160
$ws.Cells.Item(8, 3).Value = 11
$ws.Cells.Item(8, 4).Value = 26
$ws.Cells.Item(8, 5).Value = 65
$ws.Cells.Item(8, 6).Value = 73
$ws.Cells.Item(8, 7).Value = 94
$ws.Cells.Item(8, 8).Value = 138
$ws.Cells.Item(8, 9).Value = 159
$ws.Cells.Item(8, 10).Value = 193
$ws.Cells.Item(8, 11).Value = 278
$ws.Cells.Item(8, 12).Value = 359
$ws.Cells.Item(8, 13).Value = 456
$ws.Cells.Item(8, 14).Value = 548
$ws.Cells.Item(8, 15).Value = 608
$ws.Cells.Item(8, 16).Value = 704
$ws.Cells.Item(8, 17).Value = 814
$ws.Cells.Item(8, 18).Value = 2376
$ws.Cells.Item(8, 19).Value = 1
$ws.Cells.Item(8, 20).Value = 7063
